# Alumni email static data - add two new rows (rakesh_patil25, hemant_patil)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: rakesh_patil25@gmail.com / rakesh147 / teacher
$ws.Range("A8").Value = "rakesh_patil25@gmail.com"
$ws.Range("B8").Value = "rakesh147"
$ws.Range("C8").Value = "teacher"

# Row 9: hemant_patil@gmail.com / hemantpatil147 / student
$ws.Range("A9").Value = "hemant_patil@gmail.com"
$ws.Range("B9").Value = "hemantpatil147"
$ws.Range("C9").Value = "student"

# Hyperlink the two new e-mail addresses, same as the existing rows.
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:rakesh_patil25@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:hemant_patil@gmail.com")

# Match the hyperlink-style formatting already used in column A (rows 3-7).
$ws.Range("A8").Style = $ws.Range("A3").Style
$ws.Range("A9").Style = $ws.Range("A3").Style

# Leave the selection where the author left it when saving.
$ws.Range("B9").Select()
